# Add new columns I (I0) and J (IF) to the worksheet, matching the
# style used by the existing header row (B1:H1), and fill in the
# per-row values for rows 2-33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style of an existing header cell (H1) onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data values for rows 2-33 ---
$iValues = @(8,9,5,6,7,6,6,7,5,8,7,7,4,5,10,9,8,4,5,9,6,7,5,9,6,6,4,9,9,8,9,8)
$jValues = @(8,9,6,8,7,8,6,7,6,8,7,8,6,7,10,9,8,6,7,9,7,7,6,9,7,7,5,9,9,8,9,8)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}

$wb.Save()
